$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44211
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 13).Value = 1883
$ws.Cells.Item(2, 16).Value = 1883

# Row 3
$ws.Cells.Item(3, 4).Value = 44260
$ws.Cells.Item(3, 10).Value = 220
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1909
$ws.Cells.Item(3, 16).Value = 1909

# Row 4
$ws.Cells.Item(4, 4).Value = 44524
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 2000
$ws.Cells.Item(4, 13).Value = 2000
$ws.Cells.Item(4, 16).Value = 2000

# Row 5
$ws.Cells.Item(5, 4).Value = 44166
$ws.Cells.Item(5, 10).Value = 240
$ws.Cells.Item(5, 11).Value = 600
$ws.Cells.Item(5, 12).Value = 700
$ws.Cells.Item(5, 13).Value = 633
$ws.Cells.Item(5, 16).Value = 633

# Row 6
$ws.Cells.Item(6, 4).Value = 44273
$ws.Cells.Item(6, 10).Value = 140
$ws.Cells.Item(6, 13).Value = 1914
$ws.Cells.Item(6, 16).Value = 1914

# Row 7
$ws.Cells.Item(7, 4).Value = 44265
$ws.Cells.Item(7, 10).Value = 220
$ws.Cells.Item(7, 11).Value = 1800
$ws.Cells.Item(7, 13).Value = 1909
$ws.Cells.Item(7, 16).Value = 1909

# Row 8
$ws.Cells.Item(8, 4).Value = 44525
$ws.Cells.Item(8, 10).Value = 60

# Row 9
$ws.Cells.Item(9, 4).Value = 44266
$ws.Cells.Item(9, 10).Value = 150
$ws.Cells.Item(9, 11).Value = 1800
$ws.Cells.Item(9, 12).Value = 2000
$ws.Cells.Item(9, 13).Value = 1913
$ws.Cells.Item(9, 16).Value = 1913

# Row 10
$ws.Cells.Item(10, 4).Value = 44267
$ws.Cells.Item(10, 10).Value = 150
$ws.Cells.Item(10, 13).Value = 1913
$ws.Cells.Item(10, 16).Value = 1913

# Row 11
$ws.Cells.Item(11, 4).Value = 44263
$ws.Cells.Item(11, 10).Value = 140
$ws.Cells.Item(11, 13).Value = 1914
$ws.Cells.Item(11, 16).Value = 1914

# Row 12
$ws.Cells.Item(12, 4).Value = 44533
$ws.Cells.Item(12, 10).Value = 100

# Row 13
$ws.Cells.Item(13, 4).Value = 44539
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 2000
$ws.Cells.Item(13, 12).Value = 2200
$ws.Cells.Item(13, 13).Value = 2100
$ws.Cells.Item(13, 16).Value = 2100

# Row 14
$ws.Cells.Item(14, 4).Value = 44532
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 2000
$ws.Cells.Item(14, 12).Value = 2200
$ws.Cells.Item(14, 13).Value = 2100
$ws.Cells.Item(14, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(14, 16).Value = 2100

# Row 15
$ws.Cells.Item(15, 4).Value = 44271
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 1800
$ws.Cells.Item(15, 12).Value = 2000
$ws.Cells.Item(15, 13).Value = 1920
$ws.Cells.Item(15, 16).Value = 1920

# Row 16
$ws.Cells.Item(16, 4).Value = 44160
$ws.Cells.Item(16, 10).Value = 190
$ws.Cells.Item(16, 11).Value = 1300
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = 1395
$ws.Cells.Item(16, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(16, 16).Value = 930
$ws.Cells.Item(16, 17).Value = 1.5

# Row 17
$ws.Cells.Item(17, 4).Value = 44208
$ws.Cells.Item(17, 10).Value = 130
$ws.Cells.Item(17, 15).Value = "Provincia de Cautín"

# Row 18
$ws.Cells.Item(18, 4).Value = 44264
$ws.Cells.Item(18, 10).Value = 130
$ws.Cells.Item(18, 11).Value = 1800
$ws.Cells.Item(18, 13).Value = 1908
$ws.Cells.Item(18, 16).Value = 1908

# Row 19
$ws.Cells.Item(19, 4).Value = 44270
$ws.Cells.Item(19, 10).Value = 260
$ws.Cells.Item(19, 11).Value = 1800
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 1908
$ws.Cells.Item(19, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(19, 16).Value = 1908
$ws.Cells.Item(19, 17).Value = 1

# Row 20
$ws.Cells.Item(20, 4).Value = 44536
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 13).Value = 2000
$ws.Cells.Item(20, 16).Value = 2000

# Row 21
$ws.Cells.Item(21, 4).Value = 44272
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 13).Value = 1893
$ws.Cells.Item(21, 16).Value = 1893
